# Add new contact-related translation rows to the "Traductions" table.
# Cells are written in the same order as the author produced them so the
# shared-string table (insertion-ordered) comes out identical.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 64
$ws.Range("A64").Value = 'contact_msg'
$ws.Range("B64").Value = 'Merci d''avoir pris le temps de nous contacter ! '

# Row 66 (unauth_err) was filled in before row 65 (deco_msg) in the
# original edit, which is why their shared-string ids interleave.
$ws.Range("A66").Value = 'unauth_err'
$ws.Range("B66").Value = 'Vous n''êtes pas connecté(e) !'

# Row 65
$ws.Range("A65").Value = 'deco_msg'
$ws.Range("B65").Value = 'Vous avez bien été déconnecté.'

# Row 67
$ws.Range("A67").Value = 'badreq_err'
$ws.Range("B67").Value = 'Mauvaise requête.'

# Row 68
$ws.Range("A68").Value = 'co_msg'
$ws.Range("B68").Value = 'Connecté(e) !'

# Row 69
$ws.Range("A69").Value = 'contact_err'
$ws.Range("B69").Value = 'Il y a eu une erreur en traitant votre demande.'

# Row 70
$ws.Range("A70").Value = 'upload_msg'
$ws.Range("B70").Value = 'Import réussi !'

# Row 71-73 ids, then the B column values (contact_corps' value landed on
# row 74 before the row 74/75 ids were typed in).
$ws.Range("A71").Value = 'contact_nom'
$ws.Range("A72").Value = 'contact_sujet'
$ws.Range("A73").Value = 'contact_mail'

$ws.Range("B71").Value = 'Nom'
$ws.Range("B72").Value = 'Objet'
$ws.Range("B73").Value = 'Adresse mail'
$ws.Range("B74").Value = 'Message'

# Row 74-75
$ws.Range("A74").Value = 'contact_corps'
$ws.Range("A75").Value = 'contact_suppr'
$ws.Range("B75").Value = 'Supprimer'

# Resize the "Traductions" table to cover the new rows.
$lo = $ws.ListObjects.Item(1)
[void]$lo.Resize($ws.Range("A1:B75"))

# Leave the selection on the next empty row, like the author did.
[void]$ws.Range("A76").Select()
